# Applies the "Updated Okonomi spreadsheet with Person Belop" change.
#
# Summary of the functional edit (see commit message / xml diff):
#  - On "Multiple Employers": the old single "shortfall" split (row 19-26)
#    is reworked into a split between an employer shortfall adjustment and a
#    brand new "Person Belop" (person amount) adjustment. Two new rows are
#    inserted (new rows 24 and 25) and the formulas below them are rewired to
#    use three new named ranges.
#  - Three new defined names are added: AdjustmentArb, ShortfallIncomeRatioArb,
#    ShortfallIncomeRatioPerson (ScalingFactor already existed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Multiple Employers")

# ---------------------------------------------------------------------------
# 1. Insert two new rows right before the old row 24 ("Arbeidsgiverbelop").
#    This shifts the old rows 24-33 down to 26-35 and keeps their formatting
#    & (simple) formula references intact; we patch the handful of formulas
#    below that need more than a plain row-shift.
# ---------------------------------------------------------------------------
$ws.Rows(24).Insert()
$ws.Rows(24).Insert()

# ---------------------------------------------------------------------------
# 2. New defined names used by the reworked formulas.
# ---------------------------------------------------------------------------
$wb.Names.Add("AdjustmentArb", "='Multiple Employers'!`$D`$23:`$G`$23")
$wb.Names.Add("ShortfallIncomeRatioArb", "='Multiple Employers'!`$D`$20:`$G`$20")
$wb.Names.Add("ShortfallIncomeRatioPerson", "='Multiple Employers'!`$D`$24:`$G`$24")

# ---------------------------------------------------------------------------
# 3. Row 20: label the (until-now unlabeled) "shortfall income ratio" column.
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = "Shortfall income ratio - arbeidsgiver"

# ---------------------------------------------------------------------------
# 4. Row 21: relabel + rewire the formula to use the new named range.
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = "Shortfall ratio - arbeidsgiver"
$ws.Range("D21:G21").Formula = "=IF(SUM(ShortfallIncomeRatioArb)=0,0,D20/SUM(ShortfallIncomeRatioArb))"

# ---------------------------------------------------------------------------
# 5. Row 23 ("Ekstra fordeling" / "Adjustment for employers"): the shortfall
#    adjustment must never exceed the employer's own shortfall (-D19).
# ---------------------------------------------------------------------------
$ws.Range("D23:G23").Formula = "=MIN(-D19,IF(SUM(`$D`$19:`$G`$19)>=0, 0, SUM(`$D`$22:`$G`$22)*D21))"

# ---------------------------------------------------------------------------
# 6. New row 24: "Shortfall income ratio - person".
# ---------------------------------------------------------------------------
$ws.Range("B24").Value = "Shortfall income ratio - person"
$ws.Range("D24:G24").Formula = "=IF(D22 = 0, 0, D4)"

# ---------------------------------------------------------------------------
# 7. New row 25: "Shortfall ratio - person".
# ---------------------------------------------------------------------------
$ws.Range("B25").Value = "Shortfall ratio - person"
$ws.Range("D25:G25").Formula = "=IF(SUM(ShortfallIncomeRatioPerson)= 0, 0, D24/SUM(ShortfallIncomeRatioPerson))"

# ---------------------------------------------------------------------------
# 8. Row 27 ("Personbelop" / "Amount to employee"): new Person Belop amount,
#    the employer's initial amount to employee minus the employer-side
#    adjustment apportioned by the person shortfall ratio.
# ---------------------------------------------------------------------------
$ws.Range("D27:G27").Formula = "=ROUND(D22-SUM(AdjustmentArb)*D25,0)"

# ---------------------------------------------------------------------------
# 9. Row 28 ("Daily rounding error"): the rounding-error check now needs to
#    sum both the employer amount (row 26) and the person amount (row 27).
# ---------------------------------------------------------------------------
$ws.Range("C28").Formula = "=C15-SUM(D26:G27)"

# ---------------------------------------------------------------------------
# 10. Re-select/zoom to match the final view of the "Multiple Employers"
#     sheet, and set the "Terminology" sheet zoom too.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Terminology")
$ws1.Select()
$excel.ActiveWindow.Zoom = 130

$ws.Select()
$ws.Range("E11").Select()
$excel.ActiveWindow.Zoom = 130
